$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Recorded By" column (G) lists the recording entities as a comma
# separated string. Previously entries were ordered "System, <email>";
# they should now be ordered "<email>, System".
$ws.Cells.Replace(
    "System, dnasr281@gmail.com",
    "dnasr281@gmail.com, System",
    1,      # xlWhole
    1,      # xlPart (unused when LookAt is xlWhole, but keep default)
    $false,
    $false,
    $true
)
